# ADM1 DAE implementation with pH solver
# Applies:
#  - Sheet1: frozen-pane top-left cell moves from G2 to B2
#  - Sheet3: add "?" placeholders in columns K/L for several rows,
#            rewrite the Z-column rate-expression formulas (move X_xx term,
#            rename K_va/K_bu -> K_c4, rename k_ca -> k_ac), and move the
#            active selection to Z19

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: move the frozen pane's top-left visible cell from G2 to B2
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$pane = $excel.ActiveWindow.Panes.Item(4)
$pane.ScrollRow = 2
$pane.ScrollColumn = 2

# ---------------------------------------------------------------------
# Sheet3: cell content updates
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()

# New "?" placeholders (stoichiometry not yet filled in) in columns K/L
$ws3.Range("K2").Value = "?"
$ws3.Range("L2").Value = "?"

$ws3.Range("K3").Value = "?"
$ws3.Range("L3").Value = "?"

$ws3.Range("K4").Value = "?"
$ws3.Range("L4").Value = "?"

$ws3.Range("K5").Value = "?"
$ws3.Range("L5").Value = "?"

$ws3.Range("K8").Value = "?"
$ws3.Range("K9").Value = "?"
$ws3.Range("K10").Value = "?"

$ws3.Range("K14").Value = "?"
$ws3.Range("L14").Value = "?"

$ws3.Range("K15").Value = "?"
$ws3.Range("L15").Value = "?"

$ws3.Range("K16").Value = "?"
$ws3.Range("L16").Value = "?"

$ws3.Range("K17").Value = "?"
$ws3.Range("L17").Value = "?"

$ws3.Range("K18").Value = "?"
$ws3.Range("L18").Value = "?"

$ws3.Range("K19").Value = "?"
$ws3.Range("L19").Value = "?"

$ws3.Range("K20").Value = "?"
$ws3.Range("L20").Value = "?"

# Rewrite uptake-rate expressions in column Z (reorder X_xx term, rename
# K_va/K_bu to K_c4, rename k_ca to k_ac)
$ws3.Range("Z6").Value = "k_su * X_su * S_su/(K_su + S_su)"
$ws3.Range("Z7").Value = "k_aa * X_aa * S_aa/(K_aa + S_aa) "
$ws3.Range("Z8").Value = "k_fa * X_fa * S_fa/(K_fa + S_fa)"
$ws3.Range("Z9").Value = "k_c4 * X_c4 * S_va/(K_c4 + S_va) * S_va/(S_va + S_bu)"
$ws3.Range("Z10").Value = "k_c4 * X_c4 * S_bu/(K_c4 + S_bu) * S_bu/(S_va + S_bu)"
$ws3.Range("Z11").Value = "k_pro * X_pro * S_pro/(K_pro + S_pro)"
$ws3.Range("Z13").Value = "k_h2 * X_h2 * S_h2/(K_h2 + S_h2)"
$ws3.Range("Z12").Value = "k_ac * X_ac * S_ac/(K_ac + S_ac)"

# Move the active selection on Sheet3 to Z19
$ws3.Range("Z19").Select()
